$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065827479255052
$ws.Range("D2").Value = 1.070693648201878
$ws.Range("E2").Value = 1.060651646555675
$ws.Range("F2").Value = 1.078895023530078
$ws.Range("I2").Value = 1.057427365903892
$ws.Range("J2").Value = 1.070780602957231
$ws.Range("K2").Value = 1.073392804520508
$ws.Range("L2").Value = 1.063377845125482
$ws.Range("M2").Value = 1.081572497761493
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.067233380747203
$ws.Range("D3").Value = 1.071871466379218
$ws.Range("E3").Value = 1.061890197493994
$ws.Range("F3").Value = 1.080266698717642
$ws.Range("I3").Value = 1.057968143368586
$ws.Range("J3").Value = 1.071840099561328
$ws.Range("K3").Value = 1.074386139827169
$ws.Range("L3").Value = 1.064429793187528
$ws.Range("M3").Value = 1.082760793607816
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.068141393267369
$ws.Range("D4").Value = 1.072632127350014
$ws.Range("E4").Value = 1.062689604458665
$ws.Range("F4").Value = 1.081153105096139
$ws.Range("I4").Value = 1.058315908036102
$ws.Range("J4").Value = 1.072523493457386
$ws.Range("K4").Value = 1.07502685138093
$ws.Range("L4").Value = 1.065107896215758
$ws.Range("M4").Value = 1.083528000873052
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.068522721800556
$ws.Range("D5").Value = 1.072951563377022
$ws.Range("E5").Value = 1.063025198096776
$ws.Range("F5").Value = 1.081525478796881
$ws.Range("I5").Value = 1.058461595563767
$ws.Range("J5").Value = 1.072810277835612
$ws.Range("K5").Value = 1.075295722342556
$ws.Range("L5").Value = 1.065392359337207
$ws.Range("M5").Value = 1.08385013343902
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.068586725289642
$ws.Range("D6").Value = 1.073005177960183
$ws.Range("E6").Value = 1.063081517948598
$ws.Range("F6").Value = 1.081587986229358
$ws.Range("I6").Value = 1.058486027161931
$ws.Range("J6").Value = 1.072858400184984
$ws.Range("K6").Value = 1.075340838708994
$ws.Range("L6").Value = 1.065440086260244
$ws.Range("M6").Value = 1.083904197537683
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.068146490160673
$ws.Range("D7").Value = 1.072636397022394
$ws.Range("E7").Value = 1.062694090544148
$ws.Range("F7").Value = 1.081158081831319
$ws.Range("I7").Value = 1.058317856730912
$ws.Range("J7").Value = 1.072527327498916
$ws.Range("K7").Value = 1.075030445944322
$ws.Range("L7").Value = 1.065111699621175
$ws.Range("M7").Value = 1.083532306794431
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.066302966413604
$ws.Range("D8").Value = 1.071092004056841
$ws.Range("E8").Value = 1.061070642253742
$ws.Range("F8").Value = 1.079358831204239
$ws.Range("I8").Value = 1.05761057213135
$ws.Range("J8").Value = 1.07113911761919
$ws.Range("K8").Value = 1.073728932760774
$ws.Range("L8").Value = 1.063733893080732
$ws.Range("M8").Value = 1.08197444315388
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.063041109525324
$ws.Range("D9").Value = 1.068359125397496
$ws.Range("E9").Value = 1.058194185078529
$ws.Range("F9").Value = 1.076179175093205
$ws.Range("I9").Value = 1.056347607726527
$ws.Range("J9").Value = 1.06867604292126
$ws.Range("K9").Value = 1.071419635227627
$ws.Range("L9").Value = 1.061286022678951
$ws.Range("M9").Value = 1.079216023033635
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.060857132490951
$ws.Range("D10").Value = 1.06652917380132
$ws.Range("E10").Value = 1.056265580315132
$ws.Range("F10").Value = 1.074052876661716
$ws.Range("I10").Value = 1.055494250801383
$ws.Range("J10").Value = 1.067022308759848
$ws.Range("K10").Value = 1.069869127738024
$ws.Range("L10").Value = 1.059640301342321
$ws.Range("M10").Value = 1.077367807824993
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.059909116474136
$ws.Range("D11").Value = 1.065734805911138
$ws.Range("E11").Value = 1.055427785017955
$ws.Range("F11").Value = 1.073130532849826
$ws.Range("I11").Value = 1.055121994917916
$ws.Range("J11").Value = 1.066303375752429
$ws.Range("K11").Value = 1.069195068151977
$ws.Range("L11").Value = 1.058924329826847
$ws.Range("M11").Value = 1.076565233405362
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.059556621015877
$ws.Range("D12").Value = 1.065439437398486
$ws.Range("E12").Value = 1.055116177805805
$ws.Range("F12").Value = 1.072787679117253
$ws.Range("I12").Value = 1.054983306023267
$ws.Range("J12").Value = 1.066035896382689
$ws.Range("K12").Value = 1.068944283759381
$ws.Range("L12").Value = 1.058657873758373
$ws.Range("M12").Value = 1.076266771379325
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.059632248922415
$ws.Range("D13").Value = 1.065502808848505
$ws.Range("E13").Value = 1.055183037410183
$ws.Range("F13").Value = 1.072861234028026
$ws.Range("I13").Value = 1.055013074170493
$ws.Range("J13").Value = 1.066093291436907
$ws.Range("K13").Value = 1.068998096453126
$ws.Range("L13").Value = 1.058715052789407
$ws.Range("M13").Value = 1.07633080845276
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.059879986485535
$ws.Range("D14").Value = 1.065710396915108
$ws.Range("E14").Value = 1.055402035932365
$ws.Range("F14").Value = 1.073102197681397
$ws.Range("I14").Value = 1.055110539369009
$ws.Range("J14").Value = 1.06628127475003
$ws.Range("K14").Value = 1.069174346607461
$ws.Range("L14").Value = 1.058902314978353
$ws.Range("M14").Value = 1.076540569623999
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.060032577799273
$ws.Range("D15").Value = 1.065838258177894
$ws.Range("E15").Value = 1.055536913230306
$ws.Range("F15").Value = 1.073250629438296
$ws.Range("I15").Value = 1.055170535591388
$ws.Range("J15").Value = 1.066397039546239
$ws.Range("K15").Value = 1.069282885812097
$ws.Range("L15").Value = 1.059017625287171
$ws.Range("M15").Value = 1.076669763800188
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.060919999166877
$ws.Range("D16").Value = 1.0665818509763
$ws.Range("E16").Value = 1.056321124567063
$ws.Range("F16").Value = 1.074114054345926
$ws.Range("I16").Value = 1.055518898002805
$ws.Range("J16").Value = 1.067069961251117
$ws.Range("K16").Value = 1.069913805879585
$ws.Range("L16").Value = 1.059687746550444
$ws.Range("M16").Value = 1.077421023319369
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.061476022724455
$ws.Range("D17").Value = 1.067047750662324
$ws.Range("E17").Value = 1.056812312296903
$ws.Range("F17").Value = 1.074655213014781
$ws.Range("I17").Value = 1.055736678479207
$ws.Range("J17").Value = 1.067491298140651
$ws.Range("K17").Value = 1.070308843691891
$ws.Range("L17").Value = 1.060107190485043
$ws.Range("M17").Value = 1.077891651978342
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.061800116711256
$ws.Range("D18").Value = 1.06731931080893
$ws.Range("E18").Value = 1.057098554112242
$ws.Range("F18").Value = 1.074970703730838
$ws.Range("I18").Value = 1.05586344138967
$ws.Range("J18").Value = 1.067736781788538
$ws.Range("K18").Value = 1.070539004433634
$ws.Range("L18").Value = 1.060351521005054
$ws.Range("M18").Value = 1.078165941705859
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.061910586455744
$ws.Range("D19").Value = 1.067411873553642
$ws.Range("E19").Value = 1.057196111356735
$ws.Range("F19").Value = 1.075078251338502
$ws.Range("I19").Value = 1.055906619473997
$ws.Range("J19").Value = 1.067820438919537
$ws.Range("K19").Value = 1.070617439697862
$ws.Range("L19").Value = 1.060434776649931
$ws.Range("M19").Value = 1.078259430303852
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.061416390009301
$ws.Range("D20").Value = 1.06699778386449
$ws.Range("E20").Value = 1.056759639416676
$ws.Range("F20").Value = 1.074597168205323
$ws.Range("I20").Value = 1.055713340107911
$ws.Range("J20").Value = 1.06744612117609
$ws.Range("K20").Value = 1.070266486616894
$ws.Range("L20").Value = 1.060062221687683
$ws.Range("M20").Value = 1.077841180783234
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.059807043957852
$ws.Range("D21").Value = 1.065649275872284
$ws.Range("E21").Value = 1.055337557780721
$ws.Range("F21").Value = 1.073031246967475
$ws.Range("I21").Value = 1.055081849825181
$ws.Range("J21").Value = 1.066225930429374
$ws.Range("K21").Value = 1.069122456659642
$ws.Range("L21").Value = 1.058847185132019
$ws.Range("M21").Value = 1.076478809908426
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.058793096486481
$ws.Range("D22").Value = 1.064799648048384
$ws.Range("E22").Value = 1.05444104685177
$ws.Range("F22").Value = 1.072045215858027
$ws.Range("I22").Value = 1.054682395251955
$ws.Range("J22").Value = 1.065456224251664
$ws.Range("K22").Value = 1.06840079225611
$ws.Range("L22").Value = 1.058080275476423
$ws.Range("M22").Value = 1.075620203963097
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.059330810008866
$ws.Range("D23").Value = 1.065250221511545
$ws.Range("E23").Value = 1.054916533516166
$ws.Range("F23").Value = 1.072568071834037
$ws.Range("I23").Value = 1.054894383552069
$ws.Range("J23").Value = 1.065864501531899
$ws.Range("K23").Value = 1.068783586647512
$ws.Range("L23").Value = 1.058487112495504
$ws.Range("M23").Value = 1.076075561875339
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.061443336162337
$ws.Range("D24").Value = 1.067020362300863
$ws.Range("E24").Value = 1.056783440826582
$ws.Range("F24").Value = 1.074623396644631
$ws.Range("I24").Value = 1.055723886531508
$ws.Range("J24").Value = 1.067466535551781
$ws.Range("K24").Value = 1.070285626753652
$ws.Range("L24").Value = 1.060082542152964
$ws.Range("M24").Value = 1.07786398722384
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063886002974115
$ws.Range("D25").Value = 1.069067031703057
$ws.Range("E25").Value = 1.058939722532509
$ws.Range("F25").Value = 1.077002315078478
$ws.Range("I25").Value = 1.056676105916165
$ws.Range("J25").Value = 1.072018555750557
$ws.Range("K25").Value = 1.073392804520508
$ws.Range("L25").Value = 1.061921262912349
$ws.Range("M25").Value = 1.079930748593193
